$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "Valid" dates in column B (stored as date serials)
$ws.Range("B2").Value = 46063
$ws.Range("B3").Value = 44603

# Move the active selection to B4 (matches the saved sheet view state)
$ws.Range("B4").Select()
